# Add virtual generator rows (grid interconnection points) to the 'gen'
# sheet, modelled as generators with -5000..5000 MW capacity and a
# (0, 50) cost model, at buses 5, 10, 11, 26 and 27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gen")

$newGens = @(
    @{ Row = 67; Bus = 5  },
    @{ Row = 68; Bus = 10 },
    @{ Row = 69; Bus = 11 },
    @{ Row = 70; Bus = 26 },
    @{ Row = 71; Bus = 27 }
)

foreach ($g in $newGens) {
    $r = $g.Row
    $ws.Cells.Item($r, 1).Value = $g.Bus   # A: Bus ID
    $ws.Cells.Item($r, 3).Value = 5000     # C: Pg_max
    $ws.Cells.Item($r, 4).Value = -5000    # D: Pg_min
    $ws.Cells.Item($r, 5).Value = 0        # E: Qg_max
    $ws.Cells.Item($r, 6).Value = 0        # F: Qg_min
    $ws.Cells.Item($r, 7).Value = 0        # G: gen_cost_coef_0
    $ws.Cells.Item($r, 8).Value = 50       # H: gen_cost_coef_1
}

# Restore view state: the 'bus' sheet is no longer the selected tab, its
# selection and zoom are updated first, then 'gen' becomes the active
# (selected) sheet/tab with its own updated selection.

$wsBus = $wb.Worksheets.Item("bus")
$wsBus.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 130
$wsBus.Range("E36").Select() | Out-Null

$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 130
$ws.Range("K76").Select() | Out-Null
